# Refined metadata to be additional tab
#
# 1) Update the "panel_query_time" (column F) values on the existing
#    "data" sheet to reflect the new panel-query run.
# 2) Add a new "metadata" sheet (after "data") describing the panel the
#    "data" sheet was scraped from, with one header row + one data row.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# --- 1) refresh the per-row query timestamps on "data" --------------------
$data.Range("F2").Value = "2021-10-05 14:33:31.722216"
$data.Range("F3").Value = "2021-10-05 14:33:31.722224"
$data.Range("F4").Value = "2021-10-05 14:33:31.722227"
$data.Range("F5").Value = "2021-10-05 14:33:31.722230"
$data.Range("F6").Value = "2021-10-05 14:33:31.722233"
$data.Range("F7").Value = "2021-10-05 14:33:31.722236"
$data.Range("F8").Value = "2021-10-05 14:33:31.722238"
$data.Range("F9").Value = "2021-10-05 14:33:31.722241"
$data.Range("F10").Value = "2021-10-05 14:33:31.722244"
$data.Range("F11").Value = "2021-10-05 14:33:31.722246"
$data.Range("F12").Value = "2021-10-05 14:33:31.722249"
$data.Range("F13").Value = "2021-10-05 14:33:31.722251"
$data.Range("F14").Value = "2021-10-05 14:33:31.722254"
$data.Range("F15").Value = "2021-10-05 14:33:31.722256"
$data.Range("F16").Value = "2021-10-05 14:33:31.722261"
$data.Range("F17").Value = "2021-10-05 14:33:31.722290"
$data.Range("F18").Value = "2021-10-05 14:33:31.722294"
$data.Range("F19").Value = "2021-10-05 14:33:31.722297"
$data.Range("F20").Value = "2021-10-05 14:33:31.722302"
$data.Range("F21").Value = "2021-10-05 14:33:31.722305"
$data.Range("F22").Value = "2021-10-05 14:33:31.722326"
$data.Range("F23").Value = "2021-10-05 14:33:31.722360"
$data.Range("F24").Value = "2021-10-05 14:33:31.722383"
$data.Range("F25").Value = "2021-10-05 14:33:31.722387"
$data.Range("F26").Value = "2021-10-05 14:33:31.722390"
$data.Range("F27").Value = "2021-10-05 14:33:31.722393"
$data.Range("F28").Value = "2021-10-05 14:33:31.722396"
$data.Range("F29").Value = "2021-10-05 14:33:31.722399"
$data.Range("F30").Value = "2021-10-05 14:33:31.722401"
$data.Range("F31").Value = "2021-10-05 14:33:31.722404"
$data.Range("F32").Value = "2021-10-05 14:33:31.722407"
$data.Range("F33").Value = "2021-10-05 14:33:31.722410"
$data.Range("F34").Value = "2021-10-05 14:33:31.722413"
$data.Range("F35").Value = "2021-10-05 14:33:31.722415"

# --- 2) add the "metadata" sheet right after "data" ------------------------
$ws = $wb.Worksheets.Add($null, $data)
$ws.Name = "metadata"

# header row (bold / centered / bordered, matching the "data" sheet header style)
$data.Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "data_name"

$data.Range("C1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "data_id"

$data.Range("D1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "data_version"

$data.Range("E1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "data_version_created"

$data.Range("F1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "panel_query_time"

$data.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "panel_get_request"

# data row
$data.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0

$ws.Range("B2").Value = "Complement Deficiencies"
$ws.Range("C2").Value = 224
$ws.Range("D2").Value = "'0.43"
$ws.Range("E2").Value = "2021-07-22T07:27:22.165713Z"
$ws.Range("F2").Value = "2021-10-05 14:33:31.718528"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/224/?format=json"

$excel.CutCopyMode = $false
$data.Activate() | Out-Null
$data.Range("A1").Select() | Out-Null
